$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the row-4 data values per the task revision.
$ws.Range("D4").Value = 3
$ws.Range("F4").Value = 3
$ws.Range("H4").Value = 46

# Move the active selection from D5 to D4 to match the saved view state.
[void]$ws.Range("D4").Select()
